$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate Sheet1 (right after it) so the new sheet inherits the same
# sheetFormatPr / drawing relationship / styles as the rest of the workbook,
# then rename and repopulate it as "Sheet2".
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# Wipe whatever content was copied over from Sheet1.
$ws2.Range("A1:C4").ClearContents()

# Column C didn't exist on Sheet1, so stamp it with the same cell style (s="1")
# used throughout the workbook before writing any values into it.
$ws2.Range("A1").Copy()
$ws2.Range("C1:C4").PasteSpecial(-4122)

# Header row
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "passport no"
$ws2.Range("C1").Value = "address"

# Data rows
$ws2.Range("A2").Value = "anthony"
$ws2.Range("B2").Value = "a34214234"
$ws2.Range("C2").Value = "jegede street 23, opp owners occupiers"

$ws2.Range("A3").Value = "okoli"
$ws2.Range("B3").Value = "r43242424"
$ws2.Range("C3").Value = "street gonduro, appartment 34"

$ws2.Range("A4").Value = "dud"
$ws2.Range("B4").Value = "f77200033"
$ws2.Range("C4").Value = "makuerdi express, room 545"

# Column C custom width (stored width="11" chars; ColumnWidth input is offset
# by 5/6 internally by the host, so back it out to land on 11 exactly).
$ws2.Columns.Item(3).ColumnWidth = 10.166666666666666
